$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new time-tracking entries (2025-12-10 .. 2025-12-22) ---
$ws.Range("A1783").Value = "2025-12-10"
$ws.Range("B1783").Value = "12:00"
$ws.Range("C1783").Value = "18:00"
$ws.Range("D1783").Value = "6h 00m"
$ws.Range("E1783").Value = "#adoc"
$ws.Range("F1783").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1783").Value = "'True"
$ws.Range("H1783").Value = "'False"
$ws.Range("I1783").Formula = "=YEAR(A1783)"
$ws.Range("J1783").Formula = "=MONTH(A1783)"

$ws.Range("A1784").Value = "2025-12-10"
$ws.Range("B1784").Value = "19:45"
$ws.Range("C1784").Value = "20:15"
$ws.Range("D1784").Value = "0h 30m"
$ws.Range("E1784").Value = "#adoc"
$ws.Range("F1784").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1784").Value = "'True"
$ws.Range("H1784").Value = "'False"
$ws.Range("I1784").Formula = "=YEAR(A1784)"
$ws.Range("J1784").Formula = "=MONTH(A1784)"

$ws.Range("A1785").Value = "2025-12-10"
$ws.Range("B1785").Value = "22:15"
$ws.Range("C1785").Value = "00:30"
$ws.Range("D1785").Value = "2h 15m"
$ws.Range("E1785").Value = "#adoc"
$ws.Range("F1785").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1785").Value = "'True"
$ws.Range("H1785").Value = "'False"
$ws.Range("I1785").Formula = "=YEAR(A1785)"
$ws.Range("J1785").Formula = "=MONTH(A1785)"

$ws.Range("A1786").Value = "2025-12-11"
$ws.Range("B1786").Value = "14:45"
$ws.Range("C1786").Value = "19:30"
$ws.Range("D1786").Value = "4h 45m"
$ws.Range("E1786").Value = "#adoc"
$ws.Range("F1786").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1786").Value = "'True"
$ws.Range("H1786").Value = "'False"
$ws.Range("I1786").Formula = "=YEAR(A1786)"
$ws.Range("J1786").Formula = "=MONTH(A1786)"

$ws.Range("A1787").Value = "2025-12-12"
$ws.Range("B1787").Value = "13:30"
$ws.Range("C1787").Value = "18:00"
$ws.Range("D1787").Value = "3h 30m"
$ws.Range("E1787").Value = "#adoc"
$ws.Range("F1787").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1787").Value = "'True"
$ws.Range("H1787").Value = "'False"
$ws.Range("I1787").Formula = "=YEAR(A1787)"
$ws.Range("J1787").Formula = "=MONTH(A1787)"

$ws.Range("A1788").Value = "2025-12-12"
$ws.Range("B1788").Value = "23:00"
$ws.Range("C1788").Value = "00:30"
$ws.Range("D1788").Value = "1h 30m"
$ws.Range("E1788").Value = "#adoc"
$ws.Range("F1788").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1788").Value = "'True"
$ws.Range("H1788").Value = "'False"
$ws.Range("I1788").Formula = "=YEAR(A1788)"
$ws.Range("J1788").Formula = "=MONTH(A1788)"

$ws.Range("A1789").Value = "2025-12-14"
$ws.Range("B1789").Value = "14:30"
$ws.Range("C1789").Value = "20:00"
$ws.Range("D1789").Value = "5h 30m"
$ws.Range("E1789").Value = "#adoc"
$ws.Range("F1789").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1789").Value = "'True"
$ws.Range("H1789").Value = "'False"
$ws.Range("I1789").Formula = "=YEAR(A1789)"
$ws.Range("J1789").Formula = "=MONTH(A1789)"

$ws.Range("A1790").Value = "2025-12-14"
$ws.Range("B1790").Value = "20:30"
$ws.Range("C1790").Value = "22:30"
$ws.Range("D1790").Value = "2h 00m"
$ws.Range("E1790").Value = "#adoc"
$ws.Range("F1790").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1790").Value = "'True"
$ws.Range("H1790").Value = "'False"
$ws.Range("I1790").Formula = "=YEAR(A1790)"
$ws.Range("J1790").Formula = "=MONTH(A1790)"

$ws.Range("A1791").Value = "2025-12-14"
$ws.Range("B1791").Value = "22:45"
$ws.Range("C1791").Value = "23:45"
$ws.Range("D1791").Value = "1h 00m"
$ws.Range("E1791").Value = "#adoc"
$ws.Range("F1791").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1791").Value = "'True"
$ws.Range("H1791").Value = "'False"
$ws.Range("I1791").Formula = "=YEAR(A1791)"
$ws.Range("J1791").Formula = "=MONTH(A1791)"

$ws.Range("A1792").Value = "2025-12-15"
$ws.Range("B1792").Value = "13:45"
$ws.Range("C1792").Value = "19:30"
$ws.Range("D1792").Value = "5h 45m"
$ws.Range("E1792").Value = "#adoc"
$ws.Range("F1792").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1792").Value = "'True"
$ws.Range("H1792").Value = "'False"
$ws.Range("I1792").Formula = "=YEAR(A1792)"
$ws.Range("J1792").Formula = "=MONTH(A1792)"

$ws.Range("A1793").Value = "2025-12-15"
$ws.Range("B1793").Value = "20:30"
$ws.Range("C1793").Value = "21:30"
$ws.Range("D1793").Value = "1h 00m"
$ws.Range("E1793").Value = "#adoc"
$ws.Range("F1793").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1793").Value = "'True"
$ws.Range("H1793").Value = "'False"
$ws.Range("I1793").Formula = "=YEAR(A1793)"
$ws.Range("J1793").Formula = "=MONTH(A1793)"

$ws.Range("A1794").Value = "2025-12-16"
$ws.Range("B1794").Value = "12:30"
$ws.Range("C1794").Value = "19:00"
$ws.Range("D1794").Value = "6h 30m"
$ws.Range("E1794").Value = "#adoc"
$ws.Range("F1794").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1794").Value = "'True"
$ws.Range("H1794").Value = "'False"
$ws.Range("I1794").Formula = "=YEAR(A1794)"
$ws.Range("J1794").Formula = "=MONTH(A1794)"

$ws.Range("A1795").Value = "2025-12-17"
$ws.Range("B1795").Value = "12:00"
$ws.Range("C1795").Value = "20:00"
$ws.Range("D1795").Value = "8h 00m"
$ws.Range("E1795").Value = "#python"
$ws.Range("F1795").Value = "nwknowledgebase v1.0.0"
$ws.Range("G1795").Value = "'True"
$ws.Range("H1795").Value = "'True"
$ws.Range("I1795").Formula = "=YEAR(A1795)"
$ws.Range("J1795").Formula = "=MONTH(A1795)"

$ws.Range("A1796").Value = "2025-12-18"
$ws.Range("B1796").Value = "15:30"
$ws.Range("C1796").Value = "19:00"
$ws.Range("D1796").Value = "3h 30m"
$ws.Range("E1796").Value = "#python"
$ws.Range("F1796").Value = "nwreadinglist v4.4.0"
$ws.Range("G1796").Value = "'True"
$ws.Range("H1796").Value = "'False"
$ws.Range("I1796").Formula = "=YEAR(A1796)"
$ws.Range("J1796").Formula = "=MONTH(A1796)"

$ws.Range("A1797").Value = "2025-12-19"
$ws.Range("B1797").Value = "12:15"
$ws.Range("C1797").Value = "18:30"
$ws.Range("D1797").Value = "4h 15m"
$ws.Range("E1797").Value = "#python"
$ws.Range("F1797").Value = "nwreadinglist v4.4.0"
$ws.Range("G1797").Value = "'True"
$ws.Range("H1797").Value = "'False"
$ws.Range("I1797").Formula = "=YEAR(A1797)"
$ws.Range("J1797").Formula = "=MONTH(A1797)"

$ws.Range("A1798").Value = "2025-12-19"
$ws.Range("B1798").Value = "23:15"
$ws.Range("C1798").Value = "00:45"
$ws.Range("D1798").Value = "1h 30m"
$ws.Range("E1798").Value = "#python"
$ws.Range("F1798").Value = "nwreadinglist v4.4.0"
$ws.Range("G1798").Value = "'True"
$ws.Range("H1798").Value = "'False"
$ws.Range("I1798").Formula = "=YEAR(A1798)"
$ws.Range("J1798").Formula = "=MONTH(A1798)"

$ws.Range("A1799").Value = "2025-12-21"
$ws.Range("B1799").Value = "15:45"
$ws.Range("C1799").Value = "17:45"
$ws.Range("D1799").Value = "2h 00m"
$ws.Range("E1799").Value = "#python"
$ws.Range("F1799").Value = "nwreadinglist v4.4.0"
$ws.Range("G1799").Value = "'True"
$ws.Range("H1799").Value = "'False"
$ws.Range("I1799").Formula = "=YEAR(A1799)"
$ws.Range("J1799").Formula = "=MONTH(A1799)"

$ws.Range("A1800").Value = "2025-12-21"
$ws.Range("B1800").Value = "18:45"
$ws.Range("C1800").Value = "20:45"
$ws.Range("D1800").Value = "2h 00m"
$ws.Range("E1800").Value = "#python"
$ws.Range("F1800").Value = "nwreadinglist v4.4.0"
$ws.Range("G1800").Value = "'True"
$ws.Range("H1800").Value = "'False"
$ws.Range("I1800").Formula = "=YEAR(A1800)"
$ws.Range("J1800").Formula = "=MONTH(A1800)"

$ws.Range("A1801").Value = "2025-12-21"
$ws.Range("B1801").Value = "21:45"
$ws.Range("C1801").Value = "22:45"
$ws.Range("D1801").Value = "1h 00m"
$ws.Range("E1801").Value = "#python"
$ws.Range("F1801").Value = "nwreadinglist v4.4.0"
$ws.Range("G1801").Value = "'True"
$ws.Range("H1801").Value = "'False"
$ws.Range("I1801").Formula = "=YEAR(A1801)"
$ws.Range("J1801").Formula = "=MONTH(A1801)"

$ws.Range("A1802").Value = "2025-12-21"
$ws.Range("B1802").Value = "23:45"
$ws.Range("C1802").Value = "00:00"
$ws.Range("D1802").Value = "0h 15m"
$ws.Range("E1802").Value = "#python"
$ws.Range("F1802").Value = "nwreadinglist v4.4.0"
$ws.Range("G1802").Value = "'True"
$ws.Range("H1802").Value = "'False"
$ws.Range("I1802").Formula = "=YEAR(A1802)"
$ws.Range("J1802").Formula = "=MONTH(A1802)"

$ws.Range("A1803").Value = "2025-12-22"
$ws.Range("B1803").Value = "13:15"
$ws.Range("C1803").Value = "22:00"
$ws.Range("D1803").Value = "7h 00m"
$ws.Range("E1803").Value = "#python"
$ws.Range("F1803").Value = "nwreadinglist v4.4.0"
$ws.Range("G1803").Value = "'True"
$ws.Range("H1803").Value = "'True"
$ws.Range("I1803").Formula = "=YEAR(A1803)"
$ws.Range("J1803").Formula = "=MONTH(A1803)"

$ws.Range("A1804").Value = "2025-12-22"
$ws.Range("B1804").Value = "23:00"
$ws.Range("C1804").Value = "00:30"
$ws.Range("D1804").Value = "1h 30m"
$ws.Range("E1804").Value = "#python"
$ws.Range("F1804").Value = "nwreadinglist v4.4.0"
$ws.Range("G1804").Value = "'True"
$ws.Range("H1804").Value = "'True"
$ws.Range("I1804").Formula = "=YEAR(A1804)"
$ws.Range("J1804").Formula = "=MONTH(A1804)"

# --- Extend the sheet with 30 more blank (pre-formatted) rows, same as the existing tail ---
$ws.Range("A1811:J1811").Copy()
$ws.Range("A1812:J1841").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore the view: frozen header row, scrolled near the bottom, active cell G1809 ---
$ws.Range("G1809").Select()
